$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the crypto price/volume table (columns D=Price, E=Volume(1h), plus a
# couple of coin name/link swaps in B/C) to the latest scraped snapshot.
# Some "Price" values look numeric (e.g. "538.01"); a leading apostrophe is
# used for those so Excel stores them as literal text (matching the
# original inline-string cell contents) instead of auto-converting to a
# number.
$ws.Range("D2").Value = '62.995.50'
$ws.Range("E2").Value = '  -0.38%  '
$ws.Range("D3").Value = '3.072.35'
$ws.Range("E3").Value = '  +0.40%  '
$ws.Range("E4").Value = '  -0.13%  '
$ws.Range("D5").Value = '''538.01'
$ws.Range("E5").Value = '  -2.54%  '
$ws.Range("D6").Value = '''133.98'
$ws.Range("E6").Value = '  -1.39%  '
$ws.Range("D7").Value = '''1.00'
$ws.Range("E7").Value = '  -0.10%  '
$ws.Range("D8").Value = '3.064.63'
$ws.Range("E8").Value = '  +0.38%  '
$ws.Range("D9").Value = '''0.495'
$ws.Range("E9").Value = '  +1.22%  '
$ws.Range("D10").Value = '''0.155'
$ws.Range("E10").Value = '  -0.75%  '
$ws.Range("D11").Value = '''6.21'
$ws.Range("E11").Value = '  -5.64%  '
$ws.Range("D12").Value = '''0.454'
$ws.Range("E12").Value = '  +2.25%  '
$ws.Range("D13").Value = '''0.0000224'
$ws.Range("E13").Value = '  +5.08%  '
$ws.Range("D14").Value = '''34.35'
$ws.Range("E14").Value = '  -0.79%  '
$ws.Range("D15").Value = '3.565.13'
$ws.Range("E15").Value = '  -0.17%  '
$ws.Range("D16").Value = '62.968.83'
$ws.Range("E16").Value = '  -0.68%  '
$ws.Range("E17").Value = '  +0.05%  '
$ws.Range("D18").Value = '3.068.66'
$ws.Range("E18").Value = '  -0.26%  '
$ws.Range("D19").Value = '''6.65'
$ws.Range("E19").Value = '  +1.88%  '
$ws.Range("D20").Value = '''483.81'
$ws.Range("E20").Value = '  -2.59%  '
$ws.Range("D21").Value = '''13.34'
$ws.Range("E21").Value = '  -0.25%  '
$ws.Range("D22").Value = '''0.695'
$ws.Range("E22").Value = '  +0.33%  '
$ws.Range("D23").Value = '''7.14'
$ws.Range("E23").Value = '  +0.32%  '
$ws.Range("D24").Value = '''79.14'
$ws.Range("E24").Value = '  +3.41%  '
$ws.Range("D25").Value = '''12.13'
$ws.Range("E25").Value = '  -0.18%  '
$ws.Range("D26").Value = '''1.00'
$ws.Range("E26").Value = '  +0.02%  '
$ws.Range("D27").Value = '''2.70'
$ws.Range("E27").Value = '  -1.20%  '
$ws.Range("D28").Value = '''8.13'
$ws.Range("E28").Value = '  -0.09%  '
$ws.Range("D29").Value = '''0.998'
$ws.Range("E29").Value = '  -0.37%  '
$ws.Range("D30").Value = '''26.03'
$ws.Range("E30").Value = '  +0.92%  '
$ws.Range("D31").Value = '''1.87'
$ws.Range("E31").Value = '  -6.78%  '
$ws.Range("D32").Value = '''1.11'
$ws.Range("E32").Value = '  +1.12%  '
$ws.Range("B33").Value = 'Stacks'
$ws.Range("C33").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D33").Value = '''2.36'
$ws.Range("E33").Value = '  -4.73%  '
$ws.Range("B34").Value = 'OKB'
$ws.Range("C34").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D34").Value = '''56.81'
$ws.Range("E34").Value = '  -3.89%  '
$ws.Range("D35").Value = '''5.39'
$ws.Range("E35").Value = '  +6.36%  '
$ws.Range("D36").Value = '''6.00'
$ws.Range("E36").Value = '  +3.96%  '
$ws.Range("D37").Value = '''481.68'
$ws.Range("E37").Value = '  -8.78%  '
$ws.Range("B38").Value = 'VeChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D38").Value = '''0.0395'
$ws.Range("E38").Value = '  -2.84%  '
$ws.Range("B39").Value = 'Maker'
$ws.Range("C39").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D39").Value = '3.091.49'
$ws.Range("E39").Value = '  +2.12%  '
$ws.Range("D40").Value = '''0.0798'
$ws.Range("E40").Value = '  +1.95%  '
$ws.Range("D41").Value = '''0.116'
$ws.Range("E41").Value = '  -0.08%  '
$ws.Range("D42").Value = '''8.11'
$ws.Range("E42").Value = '  +1.80%  '
$ws.Range("E43").Value = '  +5.54%  '
$ws.Range("D44").Value = '''0.253'
$ws.Range("E44").Value = '  +2.37%  '
$ws.Range("E45").Value = '  +0.02%  '
$ws.Range("D46").Value = '0.0₃0542'
$ws.Range("E46").Value = '  +12.30%  '
$ws.Range("D47").Value = '''121.01'
$ws.Range("E47").Value = '  -0.15%  '
$ws.Range("D48").Value = '''2.02'
$ws.Range("E48").Value = '  -0.64%  '
$ws.Range("D49").Value = '''24.57'
$ws.Range("E49").Value = '  +3.26%  '
$ws.Range("E50").Value = '  +3.07%  '
$ws.Range("D51").Value = '''2.30'
$ws.Range("E51").Value = '  +5.93%  '
